$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Throughout the game you will be given objectives by your helmsman "
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "Your goal in the game is to navigate the ocean and complete the objectives given to you by your crew."
$find.Execute(
    $find.Text,
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    $find.Replacement.Text,
    2
)
